$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new columns P1, Q1 ---
# Copy the bold/centered/bordered formatting used by the rest of the header row (from O1)
# and apply it to the two new header cells before setting their values.
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- For data rows 2-25: swap column I<->K and column M<->O, then append P and Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value = $kVal   # I = old K
    $ws.Cells.Item($r, 11).Value = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value = 2      # P
    $ws.Cells.Item($r, 17).Value = 2      # Q
}
